# Updates cryptos list prices/volumes (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.561.45"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "1.946.99"
$ws.Range("E3").Value = "  -1.21%  "
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").Value = "'244.24"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("D7").Value = "'58.27"
$ws.Range("E7").Value = "  -5.33%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").Value = "'0.370"
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("D10").Value = "'55.85"
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("D11").Value = "'0.0843"
$ws.Range("E11").Value = "  +6.37%  "
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").Value = "'21.79"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("E14").Value = "  -3.74%  "
$ws.Range("D15").Value = "2.231.41"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").Value = "'13.63"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("D17").Value = "'5.28"
$ws.Range("E17").Value = "  -2.18%  "
$ws.Range("D18").Value = "1.942.90"
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("D19").Value = "36.442.66"
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("D20").Value = "0.0₃0875"
$ws.Range("E20").Value = "  +2.72%  "
$ws.Range("D21").Value = "'69.79"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("D22").Value = "'230.30"
$ws.Range("E22").Value = "  -3.20%  "
$ws.Range("D23").Value = "'5.03"
$ws.Range("E23").Value = "  -3.63%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("D27").Value = "'9.31"
$ws.Range("E27").Value = "  -3.76%  "
$ws.Range("D28").Value = "'162.73"
$ws.Range("E28").Value = "  +2.44%  "
$ws.Range("D29").Value = "'19.48"
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("E30").Value = "  -7.24%  "
$ws.Range("E31").Value = "  -0.71%  "
$ws.Range("D32").Value = "'1.16"
$ws.Range("E32").Value = "  +1.94%  "
$ws.Range("E33").Value = "  -3.34%  "
$ws.Range("E34").Value = "  +2.54%  "
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("D36").Value = "'6.28"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("D38").Value = "'1.77"
$ws.Range("E38").Value = "  -3.66%  "
$ws.Range("E39").Value = "  -5.07%  "
$ws.Range("D40").Value = "'3.07"
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").Value = "'0.0974"
$ws.Range("E41").Value = "  -1.42%  "
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "'2.96"
$ws.Range("E42").Value = "  +4.47%  "
$ws.Range("E43").Value = "  -3.50%  "
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").Value = "'16.11"
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("D46").Value = "1.359.57"
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("E47").Value = "  -4.52%  "
$ws.Range("D48").Value = "'88.02"
$ws.Range("E48").Value = "  -4.93%  "
$ws.Range("D49").Value = "'7.18"
$ws.Range("E49").Value = "  -4.24%  "
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("D51").Value = "'45.30"
$ws.Range("E51").Value = "  +3.65%  "
